$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Turbine MW (Max)" predictor row (row 10) entirely; rows below shift up.
$ws.Rows.Item(10).Delete()

# Updated summary statistics (rows 1-5)
$ws.Range("B1").Value = 0.7703480079554758
$ws.Range("B2").Value = 0.7312583071819396
$ws.Range("B3").Value = -0.1766480063435774
$ws.Range("B4").Value = 0.04488352281937752
$ws.Range("B5").Value = 0.11524372514674

# Updated regression table values (rows shifted up after deletion)
$ws.Range("B8").Value = 10.01964084590914
$ws.Range("C8").Value = [double]"7.50308523151725E-30"
$ws.Range("D8").Value = 259.6754510514192

$ws.Range("B9").Value = 0.002433968345001678
$ws.Range("C9").Value = 0.4850095709435331
$ws.Range("D9").Value = 3.10858594488866

$ws.Range("B10").Value = 0.000989679974744969
$ws.Range("C10").Value = 0.375805769986193
$ws.Range("D10").Value = 1.769582676427884

$ws.Range("B11").Value = -0.09345311504986409
$ws.Range("C11").Value = 0.2292920911877105
$ws.Range("D11").Value = 2.38263085379547

$ws.Range("B12").Value = -0.5184750768619794
$ws.Range("C12").Value = [double]"3.898718070655699E-06"
$ws.Range("D12").Value = 3.875755336571895

$ws.Range("B13").Value = -0.2045210812476517
$ws.Range("C13").Value = 0.0536632940471894
$ws.Range("D13").Value = 1.278706172593232

$ws.Range("B14").Value = -0.1759677541525635
$ws.Range("C14").Value = 0.2017558645641281
$ws.Range("D14").Value = 1.148672010644807

$ws.Range("B15").Value = -0.4618562439295095
$ws.Range("C15").Value = 0.02197552871617524
$ws.Range("D15").Value = 1.203169735251236

$ws.Range("B16").Value = -0.1766480063435774
$ws.Range("C16").Value = 0.0002725979274128266
$ws.Range("D16").Value = 1.483881529870047

$wb.Save()
